$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("TODOS")
$ws.Cells.Item(2, 1).Value = "16:35"
$ws.Cells.Item(2, 2).Value = "16_SANTA ANA"
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = "🚌"
$ws.Cells.Item(3, 1).Value = "16:35"
$ws.Cells.Item(3, 2).Value = "16_SANTA ANA"
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = "🚌"
$ws.Cells.Item(4, 1).Value = "16:39"
$ws.Cells.Item(4, 2).Value = "17_ROMERO"
$ws.Cells.Item(4, 3).Value = 5
$ws.Cells.Item(4, 4).Value = "📅"
$ws.Cells.Item(5, 1).Value = "16:42"
$ws.Cells.Item(5, 2).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(5, 3).Value = 8
$ws.Cells.Item(5, 4).Value = "🚌"
$ws.Cells.Item(6, 1).Value = "16:42"
$ws.Cells.Item(6, 2).Value = "225_GOMEZ"
$ws.Cells.Item(6, 3).Value = 8
$ws.Cells.Item(6, 4).Value = "📅"
$ws.Cells.Item(7, 1).Value = "16:48"
$ws.Cells.Item(7, 2).Value = "15_ABASTO"
$ws.Cells.Item(7, 3).Value = 14
$ws.Cells.Item(7, 4).Value = "🚌"
$ws.Cells.Item(8, 1).Value = "16:50"
$ws.Cells.Item(8, 2).Value = "14_ABASTO"
$ws.Cells.Item(8, 3).Value = 16
$ws.Cells.Item(8, 4).Value = "🚌"
$ws.Cells.Item(9, 1).Value = "16:51"
$ws.Cells.Item(9, 2).Value = "23_HERNANDEZ"
$ws.Cells.Item(9, 3).Value = 17
$ws.Cells.Item(9, 4).Value = "🚌"
$ws.Cells.Item(10, 1).Value = "16:52"
$ws.Cells.Item(10, 2).Value = "215B_LP-P MOR-40 Y 115"
$ws.Cells.Item(10, 3).Value = 18
$ws.Cells.Item(10, 4).Value = "🚌"
$ws.Cells.Item(11, 1).Value = "16:56"
$ws.Cells.Item(11, 2).Value = "10_OLMOS"
$ws.Cells.Item(11, 3).Value = 22
$ws.Cells.Item(11, 4).Value = "🚌"
$ws.Cells.Item(12, 1).Value = "16:56"
$ws.Cells.Item(12, 2).Value = "17_179 Y 38"
$ws.Cells.Item(12, 3).Value = 22
$ws.Cells.Item(12, 4).Value = "🚌"
$ws.Cells.Item(13, 1).Value = "17:04"
$ws.Cells.Item(13, 2).Value = "11_ETCHEVERRY"
$ws.Cells.Item(13, 3).Value = 30
$ws.Cells.Item(13, 4).Value = "🚌"
$ws.Cells.Item(14, 1).Value = "17:04"
$ws.Cells.Item(14, 2).Value = "215A_EL PATO"
$ws.Cells.Item(14, 3).Value = 30
$ws.Cells.Item(14, 4).Value = "📅"
$ws.Cells.Item(15, 1).Value = "17:04"
$ws.Cells.Item(15, 2).Value = "23_HERNANDEZ"
$ws.Cells.Item(15, 3).Value = 30
$ws.Cells.Item(15, 4).Value = "🚌"
$ws.Cells.Item(16, 1).Value = "17:09"
$ws.Cells.Item(16, 2).Value = "10_OLMOS"
$ws.Cells.Item(16, 3).Value = 35
$ws.Cells.Item(16, 4).Value = "🚌"
$ws.Cells.Item(17, 1).Value = "17:13"
$ws.Cells.Item(17, 2).Value = "215A_LA PLATA"
$ws.Cells.Item(17, 3).Value = 39
$ws.Cells.Item(17, 4).Value = "🚌"
$ws.Cells.Item(18, 1).Value = "17:16"
$ws.Cells.Item(18, 2).Value = "11_ETCHEVERRY"
$ws.Cells.Item(18, 3).Value = 42
$ws.Cells.Item(18, 4).Value = "🚌"
$ws.Cells.Item(19, 1).Value = "17:20"
$ws.Cells.Item(19, 2).Value = "26_HERNANDEZ"
$ws.Cells.Item(19, 3).Value = 46
$ws.Cells.Item(19, 4).Value = "🚌"
$ws.Cells.Item(20, 1).Value = "17:28"
$ws.Cells.Item(20, 2).Value = "14_ABASTO"
$ws.Cells.Item(20, 3).Value = 54
$ws.Cells.Item(20, 4).Value = "🚌"
$ws.Cells.Item(21, 1).Value = "17:33"
$ws.Cells.Item(21, 2).Value = "15_ABASTO"
$ws.Cells.Item(21, 3).Value = 59
$ws.Cells.Item(21, 4).Value = "🚌"
$ws.Cells.Item(22, 1).Value = "17:37"
$ws.Cells.Item(22, 2).Value = "27_EL RETIRO"
$ws.Cells.Item(22, 3).Value = 63
$ws.Cells.Item(22, 4).Value = "🚌"
$ws.Cells.Item(23, 1).Value = "17:39"
$ws.Cells.Item(23, 2).Value = "215B_EL PATO"
$ws.Cells.Item(23, 3).Value = 65
$ws.Cells.Item(23, 4).Value = "📅"
$ws.Cells.Item(24, 1).Value = "17:50"
$ws.Cells.Item(24, 2).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(24, 3).Value = 76
$ws.Cells.Item(24, 4).Value = "🚌"
$ws.Cells.Item(25, 1).Value = "17:52"
$ws.Cells.Item(25, 2).Value = "81_EL PELIGRO"
$ws.Cells.Item(25, 3).Value = 78
$ws.Cells.Item(25, 4).Value = "📅"
$ws.Cells.Item(26, 1).Value = "18:02"
$ws.Cells.Item(26, 2).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(26, 3).Value = 88
$ws.Cells.Item(26, 4).Value = "🚌"
$ws.Cells.Item(27, 1).Value = "18:03"
$ws.Cells.Item(27, 2).Value = "215C_LA PLATA"
$ws.Cells.Item(27, 3).Value = 89
$ws.Cells.Item(27, 4).Value = "🚌"
$ws.Cells.Item(28, 1).Value = "18:04"
$ws.Cells.Item(28, 2).Value = "17_ROMERO"
$ws.Cells.Item(28, 3).Value = 90
$ws.Cells.Item(28, 4).Value = "🚌"
$ws.Cells.Item(29, 1).Value = "18:20"
$ws.Cells.Item(29, 2).Value = "26_HERNANDEZ"
$ws.Cells.Item(29, 3).Value = 106
$ws.Cells.Item(29, 4).Value = "🚌"
$ws.Cells.Item(30, 1).Value = "18:27"
$ws.Cells.Item(30, 2).Value = "215C_EL PATO"
$ws.Cells.Item(30, 3).Value = 113
$ws.Cells.Item(30, 4).Value = "🚌"
$ws.Cells.Item(31, 1).Value = "18:31"
$ws.Cells.Item(31, 2).Value = "11X44_ETCHEVERRY"
$ws.Cells.Item(31, 3).Value = 117
$ws.Cells.Item(31, 4).Value = "🚌"

$ws = $wb.Worksheets.Item("215")
$ws.Cells.Item(2, 1).Value = "16:52"
$ws.Cells.Item(2, 2).Value = "215B_LP-P MOR-40 Y 115"
$ws.Cells.Item(2, 3).Value = 18
$ws.Cells.Item(2, 4).Value = "🚌"
$ws.Cells.Item(3, 1).Value = "17:04"
$ws.Cells.Item(3, 2).Value = "215A_EL PATO"
$ws.Cells.Item(3, 3).Value = 30
$ws.Cells.Item(3, 4).Value = "📅"
$ws.Cells.Item(4, 1).Value = "17:13"
$ws.Cells.Item(4, 2).Value = "215A_LA PLATA"
$ws.Cells.Item(4, 3).Value = 39
$ws.Cells.Item(4, 4).Value = "🚌"
$ws.Cells.Item(5, 1).Value = "17:39"
$ws.Cells.Item(5, 2).Value = "215B_EL PATO"
$ws.Cells.Item(5, 3).Value = 65
$ws.Cells.Item(5, 4).Value = "📅"
$ws.Cells.Item(6, 1).Value = "18:03"
$ws.Cells.Item(6, 2).Value = "215C_LA PLATA"
$ws.Cells.Item(6, 3).Value = 89
$ws.Cells.Item(6, 4).Value = "🚌"
$ws.Cells.Item(7, 1).Value = "18:27"
$ws.Cells.Item(7, 2).Value = "215C_EL PATO"
$ws.Cells.Item(7, 3).Value = 113
$ws.Cells.Item(7, 4).Value = "🚌"

$ws = $wb.Worksheets.Item("COMBINADAS")
$ws.Cells.Item(2, 1).Value = "16:35"
$ws.Cells.Item(2, 2).Value = "16_SANTA ANA"
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = "🚌"
$ws.Cells.Item(3, 1).Value = "16:35"
$ws.Cells.Item(3, 2).Value = "16_SANTA ANA"
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = "🚌"
$ws.Cells.Item(4, 1).Value = "16:39"
$ws.Cells.Item(4, 2).Value = "17_ROMERO"
$ws.Cells.Item(4, 3).Value = 5
$ws.Cells.Item(4, 4).Value = "📅"
$ws.Cells.Item(5, 1).Value = "16:42"
$ws.Cells.Item(5, 2).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(5, 3).Value = 8
$ws.Cells.Item(5, 4).Value = "🚌"
$ws.Cells.Item(6, 1).Value = "16:42"
$ws.Cells.Item(6, 2).Value = "225_GOMEZ"
$ws.Cells.Item(6, 3).Value = 8
$ws.Cells.Item(6, 4).Value = "📅"
$ws.Cells.Item(7, 1).Value = "16:48"
$ws.Cells.Item(7, 2).Value = "15_ABASTO"
$ws.Cells.Item(7, 3).Value = 14
$ws.Cells.Item(7, 4).Value = "🚌"
$ws.Cells.Item(8, 1).Value = "16:50"
$ws.Cells.Item(8, 2).Value = "14_ABASTO"
$ws.Cells.Item(8, 3).Value = 16
$ws.Cells.Item(8, 4).Value = "🚌"
$ws.Cells.Item(9, 1).Value = "16:51"
$ws.Cells.Item(9, 2).Value = "23_HERNANDEZ"
$ws.Cells.Item(9, 3).Value = 17
$ws.Cells.Item(9, 4).Value = "🚌"
$ws.Cells.Item(10, 1).Value = "16:52"
$ws.Cells.Item(10, 2).Value = "215B_LP-P MOR-40 Y 115"
$ws.Cells.Item(10, 3).Value = 18
$ws.Cells.Item(10, 4).Value = "🚌"
$ws.Cells.Item(11, 1).Value = "16:56"
$ws.Cells.Item(11, 2).Value = "10_OLMOS"
$ws.Cells.Item(11, 3).Value = 22
$ws.Cells.Item(11, 4).Value = "🚌"
$ws.Cells.Item(12, 1).Value = "16:56"
$ws.Cells.Item(12, 2).Value = "17_179 Y 38"
$ws.Cells.Item(12, 3).Value = 22
$ws.Cells.Item(12, 4).Value = "🚌"
$ws.Cells.Item(13, 1).Value = "17:04"
$ws.Cells.Item(13, 2).Value = "11_ETCHEVERRY"
$ws.Cells.Item(13, 3).Value = 30
$ws.Cells.Item(13, 4).Value = "🚌"
$ws.Cells.Item(14, 1).Value = "17:04"
$ws.Cells.Item(14, 2).Value = "215A_EL PATO"
$ws.Cells.Item(14, 3).Value = 30
$ws.Cells.Item(14, 4).Value = "📅"
$ws.Cells.Item(15, 1).Value = "17:04"
$ws.Cells.Item(15, 2).Value = "23_HERNANDEZ"
$ws.Cells.Item(15, 3).Value = 30
$ws.Cells.Item(15, 4).Value = "🚌"
$ws.Cells.Item(16, 1).Value = "17:09"
$ws.Cells.Item(16, 2).Value = "10_OLMOS"
$ws.Cells.Item(16, 3).Value = 35
$ws.Cells.Item(16, 4).Value = "🚌"
$ws.Cells.Item(17, 1).Value = "17:13"
$ws.Cells.Item(17, 2).Value = "215A_LA PLATA"
$ws.Cells.Item(17, 3).Value = 39
$ws.Cells.Item(17, 4).Value = "🚌"
$ws.Cells.Item(18, 1).Value = "17:16"
$ws.Cells.Item(18, 2).Value = "11_ETCHEVERRY"
$ws.Cells.Item(18, 3).Value = 42
$ws.Cells.Item(18, 4).Value = "🚌"
$ws.Cells.Item(19, 1).Value = "17:20"
$ws.Cells.Item(19, 2).Value = "26_HERNANDEZ"
$ws.Cells.Item(19, 3).Value = 46
$ws.Cells.Item(19, 4).Value = "🚌"
$ws.Cells.Item(20, 1).Value = "17:28"
$ws.Cells.Item(20, 2).Value = "14_ABASTO"
$ws.Cells.Item(20, 3).Value = 54
$ws.Cells.Item(20, 4).Value = "🚌"
$ws.Cells.Item(21, 1).Value = "17:33"
$ws.Cells.Item(21, 2).Value = "15_ABASTO"
$ws.Cells.Item(21, 3).Value = 59
$ws.Cells.Item(21, 4).Value = "🚌"
$ws.Cells.Item(22, 1).Value = "17:37"
$ws.Cells.Item(22, 2).Value = "27_EL RETIRO"
$ws.Cells.Item(22, 3).Value = 63
$ws.Cells.Item(22, 4).Value = "🚌"
$ws.Cells.Item(23, 1).Value = "17:39"
$ws.Cells.Item(23, 2).Value = "215B_EL PATO"
$ws.Cells.Item(23, 3).Value = 65
$ws.Cells.Item(23, 4).Value = "📅"
$ws.Cells.Item(24, 1).Value = "17:50"
$ws.Cells.Item(24, 2).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(24, 3).Value = 76
$ws.Cells.Item(24, 4).Value = "🚌"
$ws.Cells.Item(25, 1).Value = "17:52"
$ws.Cells.Item(25, 2).Value = "81_EL PELIGRO"
$ws.Cells.Item(25, 3).Value = 78
$ws.Cells.Item(25, 4).Value = "📅"
$ws.Cells.Item(26, 1).Value = "18:02"
$ws.Cells.Item(26, 2).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(26, 3).Value = 88
$ws.Cells.Item(26, 4).Value = "🚌"
$ws.Cells.Item(27, 1).Value = "18:03"
$ws.Cells.Item(27, 2).Value = "215C_LA PLATA"
$ws.Cells.Item(27, 3).Value = 89
$ws.Cells.Item(27, 4).Value = "🚌"
$ws.Cells.Item(28, 1).Value = "18:04"
$ws.Cells.Item(28, 2).Value = "17_ROMERO"
$ws.Cells.Item(28, 3).Value = 90
$ws.Cells.Item(28, 4).Value = "🚌"
$ws.Cells.Item(29, 1).Value = "18:20"
$ws.Cells.Item(29, 2).Value = "26_HERNANDEZ"
$ws.Cells.Item(29, 3).Value = 106
$ws.Cells.Item(29, 4).Value = "🚌"
$ws.Cells.Item(30, 1).Value = "18:27"
$ws.Cells.Item(30, 2).Value = "215C_EL PATO"
$ws.Cells.Item(30, 3).Value = 113
$ws.Cells.Item(30, 4).Value = "🚌"
$ws.Cells.Item(31, 1).Value = "18:31"
$ws.Cells.Item(31, 2).Value = "11X44_ETCHEVERRY"
$ws.Cells.Item(31, 3).Value = 117
$ws.Cells.Item(31, 4).Value = "🚌"

$ws = $wb.Worksheets.Item("TODOS")
$ws.Rows.Item(32).Delete()

$ws = $wb.Worksheets.Item("215")
$ws.Rows.Item(8).Delete()

$ws = $wb.Worksheets.Item("COMBINADAS")
$ws.Rows.Item(32).Delete()

Write-Host "Arribos 141 actualizados - 21"
